$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 31
$ws.Cells.Item($row, 1).Value = 9
$ws.Cells.Item($row, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item($row, 3).Value = "Metropolitana"
$ws.Cells.Item($row, 4).Value = 44461
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row - 1, 4).NumberFormat
$ws.Cells.Item($row, 5).Value = 13
$ws.Cells.Item($row, 6).Value = 100112029
$ws.Cells.Item($row, 7).Value = "Orégano"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 16
$ws.Cells.Item($row, 11).Value = 9500
$ws.Cells.Item($row, 12).Value = 10000
$ws.Cells.Item($row, 13).Value = 9750
$ws.Cells.Item($row, 14).Value = '$/docena de atados'
$ws.Cells.Item($row, 15).Value = "Región Metropolitana"
$ws.Cells.Item($row, 16).Value = 3250
$ws.Cells.Item($row, 17).Value = 3
$ws.Cells.Item($row, 18).Value = "Hortaliza"
